$d = $word.ActiveDocument

$oldText = "Relation / role type promotion. Contexts. Augmentations (of promoted players role kinds transforms)."
$newText = "Relation / role type promotion. Contexts. Augmentations (of promoted players role kinds transforms): relationship and expanded members / attributes / links / relations."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

$insertAt = 228

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Relation<Relationship<C, S, P, O>> (CSPO : Relation) Monads root hierarchy.")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Dataflow:")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Monads / Zippers (cons / graphs). Aggregation, recursion. Expressions. Signatures.")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Aggregation: nesting. Relationship C Relation holding same C context role corresponding / prefix of aggregated SPOs, same CSs for aggregated POs, etc.")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Relationship: Kinds / Roles. Aggregations: traversal /  expressions (bound functions renders CK, SK, PK, OK).")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Parent layer: current layer extension / expansion.")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Current layer: C intension, O extension.")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Next layer: current layer intension.")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Dataflow: perform augmentations on layers instantiations. Observers, observables, signatures.")

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1

$p = $d.Paragraphs($insertAt)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$insertAt = $insertAt + 1
$p = $d.Paragraphs($insertAt)
$p.Range.InsertAfter("Inferir relación dominio / rango,  alcance / campo. Describir relacion n-aria como predicados.")
